$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 10 from 45175 to 45183
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 45183
}
